$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1068.1428
$ws.Range("J28").Value = 926.3333
$ws.Range("L28").Value = 926.3333
$ws.Range("N28").Value = -1896.3333
$ws.Range("H32").Value = 3046.75
$ws.Range("J32").Value = 3476.2727
$ws.Range("L32").Value = 3476.2727
$ws.Range("N32").Value = -4128.2727
$ws.Range("H33").Value = 138.6
$ws.Range("I33").Value = 163.2
$ws.Range("J33").Value = 89.40000000000001
$ws.Range("K33").Value = 163.2
$ws.Range("L33").Value = 89.40000000000001
$ws.Range("M33").Value = 65.80000000000001
$ws.Range("N33").Value = -547.4
$ws.Range("H40").Value = 14839.565
$ws.Range("I40").Value = 22831.334
$ws.Range("J40").Value = 12018.941
$ws.Range("K40").Value = 22831.334
$ws.Range("L40").Value = 12018.941
$ws.Range("M40").Value = -22656.334
$ws.Range("N40").Value = -12368.941
$ws.Range("H46").Value = 7599.8335
$ws.Range("I46").Value = 6399.6665
$ws.Range("J46").Value = 8800
$ws.Range("K46").Value = 19198.9995
$ws.Range("L46").Value = 26400
$ws.Range("M46").Value = -19079.9995
$ws.Range("N46").Value = -26638
$ws.Range("H51").Value = 7123.304
$ws.Range("I51").Value = 4666.6665
$ws.Range("J51").Value = 8702.571
$ws.Range("K51").Value = 4666.6665
$ws.Range("L51").Value = 8702.571
$ws.Range("M51").Value = -4182.6665
$ws.Range("N51").Value = -9670.571
$ws.Range("H60").Value = 7599.8335
$ws.Range("I60").Value = 6399.6665
$ws.Range("J60").Value = 8800
$ws.Range("K60").Value = 19198.9995
$ws.Range("L60").Value = 26400
$ws.Range("M60").Value = -18714.9995
$ws.Range("N60").Value = -27368
$ws.Range("H88").Value = 8339158.5
$ws.Range("I88").Value = 33337750
$ws.Range("J88").Value = 6294.6665
$ws.Range("K88").Value = 33337750
$ws.Range("L88").Value = 6294.6665
$ws.Range("M88").Value = -33337344
$ws.Range("N88").Value = -7106.6665
$ws.Range("H91").Value = 8339158.5
$ws.Range("I91").Value = 33337750
$ws.Range("J91").Value = 6294.6665
$ws.Range("K91").Value = 33337750
$ws.Range("L91").Value = 6294.6665
$ws.Range("M91").Value = -33336346
$ws.Range("N91").Value = -9102.666499999999
$ws.Range("H135").Value = 6707.421
$ws.Range("I135").Value = 1590
$ws.Range("K135").Value = 14310
$ws.Range("M135").Value = -11775
$ws.Range("H138").Value = 6272.1113
$ws.Range("I138").Value = 2163.6428
$ws.Range("J138").Value = 7710.075
$ws.Range("K138").Value = 6490.928400000001
$ws.Range("L138").Value = 23130.225
$ws.Range("M138").Value = -1350.928400000001
$ws.Range("N138").Value = -33410.225
$ws.Range("H140").Value = 71346.82000000001
$ws.Range("J140").Value = 68963
$ws.Range("L140").Value = 68963
$ws.Range("N140").Value = -79323
$ws.Range("H141").Value = 6125.5
$ws.Range("I141").Value = 5159.727
$ws.Range("K141").Value = 15479.181
$ws.Range("M141").Value = -10299.181

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 22000
$ws.Range("J43").Value = 22000
$ws.Range("L43").Value = 22000
$ws.Range("N43").Value = -22626
$ws.Range("H45").Value = 1746.55
$ws.Range("I45").Value = 1540.3077
$ws.Range("J45").Value = 2129.5715
$ws.Range("K45").Value = 1540.3077
$ws.Range("L45").Value = 2129.5715
$ws.Range("M45").Value = -1163.3077
$ws.Range("N45").Value = -2883.5715
$ws.Range("H55").Value = 28000
$ws.Range("I55").Value = 28000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 28000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -27685
$ws.Range("N55").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1520.7693
$ws.Range("J16").Value = 1698.5
$ws.Range("L16").Value = 1698.5
$ws.Range("N16").Value = -2272.5
$ws.Range("H22").Value = 160.25
$ws.Range("J22").Value = 150.71428
$ws.Range("L22").Value = 150.71428
$ws.Range("N22").Value = -850.71428
$ws.Range("H31").Value = 3096.689
$ws.Range("I31").Value = 2910.6829
$ws.Range("J31").Value = 5003.25
$ws.Range("K31").Value = 2910.6829
$ws.Range("L31").Value = 5003.25
$ws.Range("M31").Value = -2615.6829
$ws.Range("N31").Value = -5593.25
$ws.Range("H34").Value = 3096.689
$ws.Range("I34").Value = 2910.6829
$ws.Range("J34").Value = 5003.25
$ws.Range("K34").Value = 2910.6829
$ws.Range("L34").Value = 5003.25
$ws.Range("M34").Value = -2708.6829
$ws.Range("N34").Value = -5407.25
$ws.Range("H113").Value = 1520.7693
$ws.Range("J113").Value = 1698.5
$ws.Range("L113").Value = 1698.5
$ws.Range("N113").Value = -6038.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 185.59091
$ws.Range("I38").Value = 60.25
$ws.Range("J38").Value = 336
$ws.Range("K38").Value = 180.75
$ws.Range("L38").Value = 1008
$ws.Range("M38").Value = 166.25
$ws.Range("N38").Value = -1702
$ws.Range("H118").Value = 4630.091
$ws.Range("I118").Value = 4493.1
$ws.Range("K118").Value = 13479.3
$ws.Range("M118").Value = -12236.3
$ws.Range("H131").Value = 1985
$ws.Range("I131").Value = 1804.6
$ws.Range("J131").Value = 2024.2174
$ws.Range("K131").Value = 5413.799999999999
$ws.Range("L131").Value = 6072.6522
$ws.Range("M131").Value = -373.7999999999993
$ws.Range("N131").Value = -16152.6522
$ws.Range("H136").Value = 774965.75
$ws.Range("I136").Value = 911323.2
$ws.Range("K136").Value = 2733969.6
$ws.Range("M136").Value = -2728869.6

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H113").Value = 905648.25
$ws.Range("I113").Value = 1578376.1
$ws.Range("J113").Value = 8677.777
$ws.Range("K113").Value = 1578376.1
$ws.Range("L113").Value = 8677.777
$ws.Range("M113").Value = -1576206.1
$ws.Range("N113").Value = -13017.777
$ws.Range("H132").Value = 3939.5615
$ws.Range("I132").Value = 3753.5134
$ws.Range("J132").Value = 4283.75
$ws.Range("K132").Value = 11260.5402
$ws.Range("L132").Value = 12851.25
$ws.Range("M132").Value = -8730.540199999999
$ws.Range("N132").Value = -17911.25

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2416.3572
$ws.Range("I16").Value = 1602.7
$ws.Range("J16").Value = 4450.5
$ws.Range("K16").Value = 1602.7
$ws.Range("L16").Value = 4450.5
$ws.Range("M16").Value = -1432.7
$ws.Range("N16").Value = -4790.5
$ws.Range("H22").Value = 1206.1428
$ws.Range("I22").Value = 1048.625
$ws.Range("J22").Value = 1303.0769
$ws.Range("K22").Value = 1048.625
$ws.Range("L22").Value = 1303.0769
$ws.Range("M22").Value = -753.625
$ws.Range("N22").Value = -1893.0769
$ws.Range("H27").Value = 1206.1428
$ws.Range("I27").Value = 1048.625
$ws.Range("J27").Value = 1303.0769
$ws.Range("K27").Value = 1048.625
$ws.Range("L27").Value = 1303.0769
$ws.Range("M27").Value = -941.625
$ws.Range("N27").Value = -1517.0769
$ws.Range("H100").Value = 1414.3334
$ws.Range("I100").Value = 1556.4
$ws.Range("K100").Value = 1556.4
$ws.Range("M100").Value = -1015.4
$ws.Range("H132").Value = 4925.0806
$ws.Range("I132").Value = 4263.2095
$ws.Range("K132").Value = 12789.6285
$ws.Range("M132").Value = -10259.6285
$ws.Range("H136").Value = 5358.511
$ws.Range("I136").Value = 4924.2705
$ws.Range("K136").Value = 14772.8115
$ws.Range("M136").Value = -12222.8115

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13230825
$ws.Range("I132").Value = 1737617.4
$ws.Range("K132").Value = 5212852.199999999
$ws.Range("M132").Value = -5210322.199999999
$ws.Range("H139").Value = 58678.25
$ws.Range("J139").Value = 57357.5
$ws.Range("L139").Value = 57357.5
$ws.Range("N139").Value = -67637.5
